$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, reusing the same formatting as the other
# header cells (e.g. G1 "sum") rather than building a brand new style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the corresponding data value in H2 (unformatted, like G2)
$ws.Range("H2").Value = 0
